# Update invoice data rows 2-21 with new values (cleanup condition change)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row=2;  A="Invoice 269.pdf"; B="8/17/2016";  C=5004; D="Ricardo Noice";      E=400;  F=600; G=1000 },
    @{ Row=3;  A="Invoice 270.pdf"; B="5/26/2018";  C=7953; D="Sandye Wankel";      E=600;  F=450; G=1050 },
    @{ Row=4;  A="Invoice 271.pdf"; B="6/21/2016";  C=5556; D="Hagan Ledwitch";     E=1800; F=675; G=2475 },
    @{ Row=5;  A="Invoice 272.pdf"; B="7/14/2017";  C=7232; D="Gale Yelland";       E=400;  F=750; G=1150 },
    @{ Row=6;  A="Invoice 273.pdf"; B="12/7/2016";  C=4425; D="Marcelo Dewan";      E=1000; F=300; G=1300 },
    @{ Row=7;  A="Invoice 274.pdf"; B="10/28/2017"; C=7181; D="Lena Hughes";        E=800;  F=675; G=1475 },
    @{ Row=8;  A="Invoice 275.pdf"; B="3/17/2018";  C=5242; D="Dulci Scay";         E=200;  F=675; G=875  },
    @{ Row=9;  A="Invoice 276.pdf"; B="9/27/2016";  C=7397; D="Karlen Gavin";       E=600;  F=300; G=900  },
    @{ Row=10; A="Invoice 277.pdf"; B="1/3/2016";   C=3823; D="Marjie De";          E=600;  F=450; G=1050 },
    @{ Row=11; A="Invoice 278.pdf"; B="6/17/2018";  C=6235; D="Sianna Lavrinov";    E=1600; F=525; G=2125 },
    @{ Row=12; A="Invoice 279.pdf"; B="12/11/2018"; C=1375; D="Mycah McIver";       E=1400; F=375; G=1775 },
    @{ Row=13; A="Invoice 280.pdf"; B="8/17/2018";  C=4072; D="Manon Yele";         E=200;  F=300; G=500  },
    @{ Row=14; A="Invoice 281.pdf"; B="3/24/2016";  C=6656; D="Priscella Paireman"; E=1800; F=675; G=2475 },
    @{ Row=15; A="Invoice 282.pdf"; B="11/2/2018";  C=1191; D="Stormie Nazair";     E=1800; F=450; G=2250 },
    @{ Row=16; A="Invoice 283.pdf"; B="6/19/2017";  C=6376; D="Tiebout Gatenby";    E=1000; F=375; G=1375 },
    @{ Row=17; A="Invoice 284.pdf"; B="12/22/2018"; C=2264; D="Bea Dyte";           E=200;  F=600; G=800  },
    @{ Row=18; A="Invoice 285.pdf"; B="7/20/2017";  C=5789; D="Sancho Roxbee";      E=400;  F=375; G=775  },
    @{ Row=19; A="Invoice 286.pdf"; B="12/17/2016"; C=5873; D="Wildon Brampton";    E=200;  F=600; G=800  },
    @{ Row=20; A="Invoice 287.pdf"; B="7/4/2018";   C=8257; D="Marketa Soeiro";     E=1000; F=225; G=1225 },
    @{ Row=21; A="Invoice 288.pdf"; B="3/18/2018";  C=8084; D="Cassius Cassley";    E=1800; F=750; G=2550 }
)

foreach ($r in $data) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}
